# Add New TC for Sim ATM
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the workbook/file rename (BD -> LD)
$ws.Name = "Simatm26_1_1_3_LD"

# Add the new "run" flag value in column A of the new test-case row
$ws.Range("A2").Value = "run"

# Move the active selection to reflect where the author left off editing
$ws.Range("C29").Select() | Out-Null
